# Append: 2026-01-24 18:35 JST
# Update the "取得日時" (acquired datetime) column A for the existing
# data rows (2-7) on the "ランサーズ" sheet from the previous timestamp
# to the new run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2026-01-24 18:27:25"
$newTimestamp = "2026-01-24 18:35:17"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 7 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
